# Refresh cached market-board figures on the Leve profit sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N)
# per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2670.625
$ws.Range("J17").Value = 2360.8333
$ws.Range("L17").Value = 7082.499899999999
$ws.Range("N17").Value = -7418.499899999999
$ws.Range("H40").Value = 7074.5
$ws.Range("I40").Value = 6099.6665
$ws.Range("J40").Value = 9999
$ws.Range("K40").Value = 6099.6665
$ws.Range("L40").Value = 9999
$ws.Range("M40").Value = -5924.6665
$ws.Range("N40").Value = -10349
$ws.Range("H138").Value = 3508.2927
$ws.Range("I138").Value = 3061.2917
$ws.Range("K138").Value = 9183.875100000001
$ws.Range("M138").Value = -4043.875100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4233.391
$ws.Range("I32").Value = 4233.391
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4233.391
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3946.391
$ws.Range("N32").ClearContents()
$ws.Range("H132").Value = 3040.7446
$ws.Range("I132").Value = 2800.5264
$ws.Range("J132").Value = 4055
$ws.Range("K132").Value = 8401.5792
$ws.Range("L132").Value = 12165
$ws.Range("M132").Value = -5871.5792
$ws.Range("N132").Value = -17225

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 539.2222
$ws.Range("I22").Value = 387.66666
$ws.Range("K22").Value = 387.66666
$ws.Range("M22").Value = -214.66666
$ws.Range("H99").Value = 4191.7407
$ws.Range("I99").Value = 2750.5557
$ws.Range("J99").Value = 7074.1113
$ws.Range("K99").Value = 2750.5557
$ws.Range("L99").Value = 7074.1113
$ws.Range("M99").Value = -1252.5557
$ws.Range("N99").Value = -10070.1113
$ws.Range("H105").Value = 11820606
$ws.Range("I105").Value = 626802.25
$ws.Range("J105").Value = 41670748
$ws.Range("K105").Value = 626802.25
$ws.Range("L105").Value = 41670748
$ws.Range("M105").Value = -625055.25
$ws.Range("N105").Value = -41674242
$ws.Range("H107").Value = 3497935.5
$ws.Range("I107").Value = 5495720
$ws.Range("J107").Value = 1813
$ws.Range("K107").Value = 5495720
$ws.Range("L107").Value = 1813
$ws.Range("M107").Value = -5493800
$ws.Range("N107").Value = -5653
$ws.Range("H134").Value = 2967.8293
$ws.Range("I134").Value = 2657.8057
$ws.Range("K134").Value = 7973.4171
$ws.Range("M134").Value = -5438.4171
$ws.Range("H138").Value = 47515.47
$ws.Range("J138").Value = 47515.47
$ws.Range("L138").Value = 47515.47
$ws.Range("N138").Value = -57795.47
$ws.Range("H140").Value = 60550
$ws.Range("J140").Value = 60550
$ws.Range("L140").Value = 60550
$ws.Range("N140").Value = -70910

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3475.2246
$ws.Range("I31").Value = 2606.775
$ws.Range("K31").Value = 2606.775
$ws.Range("M31").Value = -2311.775
$ws.Range("H34").Value = 3475.2246
$ws.Range("I34").Value = 2606.775
$ws.Range("K34").Value = 2606.775
$ws.Range("M34").Value = -2404.775
$ws.Range("H58").Value = 2556.125
$ws.Range("I58").Value = 1374.75
$ws.Range("K58").Value = 1374.75
$ws.Range("M58").Value = -1171.75
$ws.Range("H105").Value = 1961.4
$ws.Range("I105").Value = 1313.75
$ws.Range("J105").Value = 2393.1667
$ws.Range("K105").Value = 1313.75
$ws.Range("L105").Value = 2393.1667
$ws.Range("M105").Value = 433.25
$ws.Range("N105").Value = -5887.1667
$ws.Range("H132").Value = 26317582
$ws.Range("I132").Value = 41667724
$ws.Range("K132").Value = 125003172
$ws.Range("M132").Value = -125000642
$ws.Range("H136").Value = 2556.125
$ws.Range("I136").Value = 1374.75
$ws.Range("K136").Value = 4124.25
$ws.Range("M136").Value = -1574.25
$ws.Range("H141").Value = 609331.1
$ws.Range("J141").Value = 609331.1
$ws.Range("L141").Value = 609331.1
$ws.Range("N141").Value = -619691.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1721.5
$ws.Range("I46").Value = 1999
$ws.Range("J46").Value = 1629
$ws.Range("K46").Value = 5997
$ws.Range("L46").Value = 4887
$ws.Range("M46").Value = -5906
$ws.Range("N46").Value = -5069
$ws.Range("H60").Value = 4280.467
$ws.Range("J60").Value = 4690
$ws.Range("L60").Value = 14070
$ws.Range("N60").Value = -14572
$ws.Range("H107").Value = 254.07143
$ws.Range("J107").Value = 254.07143
$ws.Range("L107").Value = 762.21429
$ws.Range("N107").Value = -4602.21429
$ws.Range("H132").Value = 1774.75
$ws.Range("I132").Value = 1450.5
$ws.Range("J132").Value = 2099
$ws.Range("K132").Value = 13054.5
$ws.Range("L132").Value = 18891
$ws.Range("M132").Value = -10524.5
$ws.Range("N132").Value = -23951

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 228721.11
$ws.Range("I70").Value = 503373.5
$ws.Range("J70").Value = 8999.200000000001
$ws.Range("K70").Value = 503373.5
$ws.Range("L70").Value = 8999.200000000001
$ws.Range("M70").Value = -503103.5
$ws.Range("N70").Value = -9539.200000000001
$ws.Range("H73").Value = 228721.11
$ws.Range("I73").Value = 503373.5
$ws.Range("J73").Value = 8999.200000000001
$ws.Range("K73").Value = 503373.5
$ws.Range("L73").Value = 8999.200000000001
$ws.Range("M73").Value = -502437.5
$ws.Range("N73").Value = -10871.2
$ws.Range("H80").Value = 76925880
$ws.Range("I80").Value = 111112856
$ws.Range("K80").Value = 111112856
$ws.Range("M80").Value = -111111858
$ws.Range("H83").Value = 76925880
$ws.Range("I83").Value = 111112856
$ws.Range("K83").Value = 555564280
$ws.Range("M83").Value = -555559288
$ws.Range("H107").Value = 9457
$ws.Range("I107").Value = 2230.5715
$ws.Range("K107").Value = 2230.5715
$ws.Range("M107").Value = -310.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4198.9165
$ws.Range("I40").Value = 3505.6956
$ws.Range("K40").Value = 3505.6956
$ws.Range("M40").Value = -3369.6956
$ws.Range("H46").Value = 2633.5
$ws.Range("J46").Value = 2799.125
$ws.Range("L46").Value = 2799.125
$ws.Range("N46").Value = -3175.125
$ws.Range("H61").Value = 27568
$ws.Range("I61").Value = 4014.7144
$ws.Range("J61").Value = 110004.5
$ws.Range("K61").Value = 4014.7144
$ws.Range("L61").Value = 110004.5
$ws.Range("M61").Value = -3812.7144
$ws.Range("N61").Value = -110408.5
$ws.Range("H100").Value = 3398.6365
$ws.Range("I100").Value = 2527.1428
$ws.Range("J100").Value = 4923.75
$ws.Range("K100").Value = 2527.1428
$ws.Range("L100").Value = 4923.75
$ws.Range("M100").Value = -1986.1428
$ws.Range("N100").Value = -6005.75
$ws.Range("H113").Value = 27568
$ws.Range("I113").Value = 4014.7144
$ws.Range("J113").Value = 110004.5
$ws.Range("K113").Value = 4014.7144
$ws.Range("L113").Value = 110004.5
$ws.Range("M113").Value = -1844.7144
$ws.Range("N113").Value = -114344.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 166667630
$ws.Range("I100").Value = 100000770
$ws.Range("J100").Value = 250001230
$ws.Range("K100").Value = 200001540
$ws.Range("L100").Value = 500002460
$ws.Range("M100").Value = -200000999
$ws.Range("N100").Value = -500003542

Write-Output "Leve profit refresh applied: 180 cells updated, 1 cleared"
